# Insert a new data row right before the existing row 452, shifting
# the existing rows 452:493 down to 453:494 (dimension grows to A1:R494).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("452:452").Insert()

# Populate the newly inserted row 452 with the new record's values.
$ws.Range("A452").Value = 4
$ws.Range("B452").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C452").Value = "Los Lagos"
$ws.Range("D452").Value = 45166
$ws.Range("E452").Value = 10
$ws.Range("F452").Value = 100112043
$ws.Range("G452").Value = "Pepino ensalada"
$ws.Range("H452").Value = "Sin especificar"
$ws.Range("I452").Value = "Primera"
$ws.Range("J452").Value = 120
$ws.Range("K452").Value = 14000
$ws.Range("L452").Value = 14000
$ws.Range("M452").Value = 14000
$ws.Range("N452").Value = "$/caja 60 unidades"
$ws.Range("O452").Value = "Región de Arica y Parinacota"
$ws.Range("P452").Value = 233
$ws.Range("Q452").Value = 60
$ws.Range("R452").Value = "Hortaliza"
